$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44; this shifts the existing rows 44-67 down to 45-68,
# preserving all of their data/formatting (matches the diff which re-numbers rows 44-67 -> 45-68).
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new weekly price record.
$ws.Cells.Item(44, 1).Value  = 7
$ws.Cells.Item(44, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(44, 3).Value  = "Ñuble"
$ws.Cells.Item(44, 4).Value  = 44942
$ws.Cells.Item(44, 5).Value  = 16
$ws.Cells.Item(44, 6).Value  = "Fruta"
$ws.Cells.Item(44, 7).Value  = 100103
$ws.Cells.Item(44, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(44, 9).Value  = 100103002
$ws.Cells.Item(44, 10).Value = "Ciruela"
$ws.Cells.Item(44, 11).Value = "Black Amber"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 60
$ws.Cells.Item(44, 14).Value = 15000
$ws.Cells.Item(44, 15).Value = 16000
$ws.Cells.Item(44, 16).Value = 15500
$ws.Cells.Item(44, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(44, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(44, 19).Value = 861
$ws.Cells.Item(44, 20).Value = 18
